$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new "Date of Birth" column
$ws.Range("R1").Value = "Birth Dt (O)"

# Fill the data rows (2-93) with the sample date of birth value, stored
# as text the same way the original export stored its other date-like
# text columns (shared string), matching the source extract behaviour.
$ws.Range("R2:R93").Value = "19/03/1985"

# Build a date-formatted, left-aligned style (on a scratch cell outside
# the used range) and apply it in one shot to the whole R2:R93 block so
# every cell shares a single style entry instead of one-off styles.
$scratch = $ws.Range("T1")
$scratch.HorizontalAlignment = -4131
$scratch.NumberFormat = "mm-dd-yy"
$ws.Range("R2:R93").Style = $scratch.Style
$scratch.Clear()
